# Dispatch Summary — role-based access controls + refreshed dispatch rows.
# The sheet now carries 13 dispatch rows (rows 2-14) instead of 6, each row
# keeping Order ID (A), Customer (B), Product (C), Ordered/Dispatched Qty
# (D/E), Price per Unit (F), Unit Type (G), Created At (H), Dispatched At
# (I), Salesperson (J) and Dispatched By (K, new role-based dispatcher).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  A=13; B='madhu';        C='black salt';    D='100 KG';  E='100 KG';  F='35.0 INR'; G='Per KG'; H='29-04-2025 02:21 PM'; I='29-04-2025 02:45 PM'; J='vishal.sharma'; K='ajay.sharma' }
    @{ Row=3;  A=12; B='aakarsh';      C='sendha namak';  D='100 KG';  E='100 KG';  F='22.0 INR'; G='Per KG'; H='29-04-2025 02:10 PM'; I='29-04-2025 02:10 PM'; J='admin';         K='admin' }
    @{ Row=4;  A=2;  B='praveen & co'; C='MTV Rock salt'; D='10 kg';   E='10 kg';   F='None None';            H='28-04-2025 04:26 PM'; I='29-04-2025 10:10 AM'; J='admin';         K='admin' }
    @{ Row=5;  A=8;                                        D='1 KG';   E='1 KG';    F='0.0 INR';  G='Per KG'; H='28-04-2025 07:27 PM'; I='29-04-2025 10:10 AM'; J='admin';         K='admin' }
    @{ Row=6;  A=10;                                        D='1 KG';   E='1 KG';    F='0.0 INR';  G='Per KG'; H='29-04-2025 09:54 AM'; I='29-04-2025 10:10 AM'; J='admin';         K='admin' }
    @{ Row=7;  A=11; B='Madhu123';     C='rock salt';     D='100 KG';  E='100 KG';  F='22.0 INR'; G='Per KG'; H='29-04-2025 09:54 AM'; I='29-04-2025 10:09 AM'; J='admin';         K='admin' }
    @{ Row=8;  A=9;  B='pkc';          C='cheetah salt';  D='1000 KG'; E='1000 KG'; F='11.0 INR'; G='Per KG'; H='29-04-2025 08:35 AM'; I='29-04-2025 08:35 AM'; J='admin';         K='admin' }
    @{ Row=9;  A=6;  B='amit';         C='Virat rock Salt'; D='1000 KG'; E='1000 KG'; F='10.0 INR'; G='Per KG'; H='28-04-2025 06:37 PM'; I='28-04-2025 06:48 PM'; J='admin';      K='admin' }
    @{ Row=10; A=7;  B='Amit Jawla';   C='Chheta';        D='1000 KG'; E='1000 KG'; F='11.0 INR'; G='Per KG'; H='28-04-2025 06:45 PM'; I='28-04-2025 06:48 PM'; J='admin';         K='admin' }
    @{ Row=11; A=3;  B='praveen & co'; C='MTV Rock salt'; D='20 kg';   E='20 kg';   F='None None';            H='28-04-2025 04:29 PM'; I='28-04-2025 05:46 PM'; J='admin';         K='admin' }
    @{ Row=12; A=4;  B='praveen & co'; C='MTV Rock salt'; D='30 kg';   E='30 kg';   F='None None';            H='28-04-2025 04:35 PM'; I='28-04-2025 05:46 PM'; J='admin';         K='admin' }
    @{ Row=13; A=5;  B='praveen & co'; C='MTV Rock salt'; D='100 kg';  E='100 kg';  F='None None';            H='28-04-2025 04:43 PM'; I='28-04-2025 04:43 PM'; J='admin' }
    @{ Row=14; A=1;  B='praveen & co'; C='MTV Rock salt'; D='10 kg';   E='10 kg';   F='None None';            H='28-04-2025 03:58 PM'; I='28-04-2025 04:36 PM'; J='admin' }
)

# Columns that must be blank on rows 5 and 6 (no salesperson assigned yet).
$blankCols = @('B', 'C')
$allCols = @('A', 'B', 'C', 'D', 'E', 'F', 'G', 'H', 'I', 'J', 'K')

foreach ($r in $rows) {
    $rowNum = $r.Row
    foreach ($col in $allCols) {
        $cell = $ws.Range("$col$rowNum")
        if ($r.ContainsKey($col)) {
            $cell.Value = $r[$col]
        } elseif ($blankCols -contains $col) {
            $cell.ClearContents()
        }
    }
}
